$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.104.33'
$ws.Range('E2').Value = '  +6.97%  '
$ws.Range('D3').Value = '3.685.36'
$ws.Range('E3').Value = '  +19.79%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.16'
$ws.Range('E5').Value = '  +4.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '185.04'
$ws.Range('E6').Value = '  +8.51%  '
$ws.Range('D7').Value = '3.681.65'
$ws.Range('E7').Value = '  +19.74%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.537'
$ws.Range('E9').Value = '  +5.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.163'
$ws.Range('E10').Value = '  +9.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.57'
$ws.Range('E11').Value = '  +4.76%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.500'
$ws.Range('E12').Value = '  +7.52%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '39.88'
$ws.Range('E13').Value = '  +12.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000253'
$ws.Range('E14').Value = '  +6.77%  '
$ws.Range('D15').Value = '4.299.84'
$ws.Range('E15').Value = '  +19.85%  '
$ws.Range('D16').Value = '3.687.58'
$ws.Range('E16').Value = '  +19.91%  '
$ws.Range('D17').Value = '71.082.73'
$ws.Range('E17').Value = '  +7.03%  '
$ws.Range('E18').Value = '  +2.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.52'
$ws.Range('E19').Value = '  +8.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.02'
$ws.Range('E20').Value = '  +1.96%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '521.28'
$ws.Range('E21').Value = '  +8.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.24'
$ws.Range('E22').Value = '  +17.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.745'
$ws.Range('E23').Value = '  +9.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '87.97'
$ws.Range('E24').Value = '  +6.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.49'
$ws.Range('E25').Value = '  +7.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.41'
$ws.Range('E26').Value = '  +9.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.82'
$ws.Range('E27').Value = '  +8.68%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.55'
$ws.Range('E29').Value = '  +14.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.14'
$ws.Range('E30').Value = '  +3.15%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '31.90'
$ws.Range('E31').Value = '  +14.82%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.77'
$ws.Range('E32').Value = '  +7.54%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0000110'
$ws.Range('E33').Value = '  +20.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.116'
$ws.Range('E34').Value = '  +5.19%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.20'
$ws.Range('E36').Value = '  +11.79%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.01'
$ws.Range('E37').Value = '  +8.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.343'
$ws.Range('E38').Value = '  +12.80%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.12'
$ws.Range('E39').Value = '  +9.10%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '50.96'
$ws.Range('E40').Value = '  +4.40%  '
$ws.Range('B41').Value = 'Arweave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '46.48'
$ws.Range('E41').Value = '  -3.97%  '
$ws.Range('E42').Value = '  +5.00%  '
$ws.Range('D43').Value = '3.178.98'
$ws.Range('E43').Value = '  +14.98%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.89'
$ws.Range('E44').Value = '  +8.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.78'
$ws.Range('E45').Value = '  +9.62%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '28.64'
$ws.Range('E46').Value = '  +18.71%  '
$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '396.05'
$ws.Range('E47').Value = '  +9.03%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0368'
$ws.Range('E48').Value = '  +7.97%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '135.33'
$ws.Range('E50').Value = '  +0.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.47'
$ws.Range('E51').Value = '  +14.26%  '
